$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 62, shifting existing rows 62:79 down to 63:80
# (mirrors a new weekly record being prepended to this price-history block)
$ws.Rows.Item(62).Insert()

# Populate the newly inserted row 62 with the new record's values
$ws.Range("A62").Value = 9
$ws.Range("B62").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C62").Value = "Metropolitana"
$ws.Range("D62").Value = 45135
$ws.Range("E62").Value = 13
$ws.Range("F62").Value = 100112010
$ws.Range("G62").Value = "Achicoria"
$ws.Range("H62").Value = "Sin especificar"
$ws.Range("I62").Value = "Primera"
$ws.Range("J62").Value = 90
$ws.Range("K62").Value = 7000
$ws.Range("L62").Value = 7000
$ws.Range("M62").Value = 7000
$ws.Range("N62").Value = "$/caja 16 unidades"
$ws.Range("O62").Value = "Provincia de Quillota"
$ws.Range("P62").Value = 438
$ws.Range("Q62").Value = 16
$ws.Range("R62").Value = "Hortaliza"

# Ensure the date cell keeps the workbook's date number format (style index 2)
$ws.Range("D62").NumberFormat = $ws.Range("D63").NumberFormat
